$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")

# ALC row 40
$ws_ALC.Range("H40").Value = 4710
$ws_ALC.Range("I40").Value = 3975.25
$ws_ALC.Range("J40").Value = 5444.75
$ws_ALC.Range("K40").Value = 3975.25
$ws_ALC.Range("L40").Value = 5444.75
$ws_ALC.Range("M40").Value = -3800.25
$ws_ALC.Range("N40").Value = -5794.75

# ALC row 51
$ws_ALC.Range("H51").Value = 17624.375
$ws_ALC.Range("I51").Value = 23200
$ws_ALC.Range("J51").Value = 8331.666999999999
$ws_ALC.Range("K51").Value = 23200
$ws_ALC.Range("L51").Value = 8331.666999999999
$ws_ALC.Range("M51").Value = -22716
$ws_ALC.Range("N51").Value = -9299.666999999999

# ALC row 112
$ws_ALC.Range("H112").Value = 2898.75
$ws_ALC.Range("I112").Value = 0
$ws_ALC.Range("J112").Value = 2898.75
$ws_ALC.Range("K112").Value = 0
$ws_ALC.Range("L112").Value = 8696.25
$ws_ALC.Range("M112").Value = ""
$ws_ALC.Range("N112").Value = -10912.25

# ALC row 137
$ws_ALC.Range("H137").Value = 1762.4193
$ws_ALC.Range("I137").Value = 1167.7333
$ws_ALC.Range("K137").Value = 3503.199900000001
$ws_ALC.Range("M137").Value = -953.1999000000005

# ARM row 32
$ws_ARM.Range("H32").Value = 33758.35
$ws_ARM.Range("I32").Value = 35009.16
$ws_ARM.Range("K32").Value = 35009.16
$ws_ARM.Range("M32").Value = -34722.16

# ARM row 61
$ws_ARM.Range("H61").Value = 2076.3333
$ws_ARM.Range("I61").Value = 2129.45
$ws_ARM.Range("K61").Value = 2129.45
$ws_ARM.Range("M61").Value = -1917.45

# ARM row 110
$ws_ARM.Range("H110").Value = 2273.8462
$ws_ARM.Range("I110").Value = 2046.6666
$ws_ARM.Range("K110").Value = 2046.6666
$ws_ARM.Range("M110").Value = -1.666600000000017

# ARM row 114
$ws_ARM.Range("H114").Value = 56074.5
$ws_ARM.Range("J114").Value = 56074.5
$ws_ARM.Range("L114").Value = 56074.5
$ws_ARM.Range("N114").Value = -64752.5

# ARM row 132
$ws_ARM.Range("H132").Value = 26088.781
$ws_ARM.Range("I132").Value = 26088.781
$ws_ARM.Range("J132").Value = 0
$ws_ARM.Range("K132").Value = 78266.34299999999
$ws_ARM.Range("L132").Value = 0
$ws_ARM.Range("M132").Value = -75736.34299999999
$ws_ARM.Range("N132").Value = ""

# ARM row 136
$ws_ARM.Range("H136").Value = 2076.3333
$ws_ARM.Range("I136").Value = 2129.45
$ws_ARM.Range("K136").Value = 6388.349999999999
$ws_ARM.Range("M136").Value = -3838.349999999999

# BSM row 81
$ws_BSM.Range("H81").Value = 76995
$ws_BSM.Range("J81").Value = 76995
$ws_BSM.Range("L81").Value = 76995
$ws_BSM.Range("N81").Value = -79117

# BSM row 84
$ws_BSM.Range("H84").Value = 76995
$ws_BSM.Range("J84").Value = 76995
$ws_BSM.Range("L84").Value = 230985
$ws_BSM.Range("N84").Value = -241593

# BSM row 86
$ws_BSM.Range("H86").Value = 3508.6924
$ws_BSM.Range("I86").Value = 2077.111
$ws_BSM.Range("J86").Value = 4266.5884
$ws_BSM.Range("K86").Value = 2077.111
$ws_BSM.Range("L86").Value = 4266.5884
$ws_BSM.Range("M86").Value = -954.1109999999999
$ws_BSM.Range("N86").Value = -6512.5884

# BSM row 89
$ws_BSM.Range("H89").Value = 3508.6924
$ws_BSM.Range("I89").Value = 2077.111
$ws_BSM.Range("J89").Value = 4266.5884
$ws_BSM.Range("K89").Value = 10385.555
$ws_BSM.Range("L89").Value = 21332.942
$ws_BSM.Range("M89").Value = -4769.555
$ws_BSM.Range("N89").Value = -32564.942

# BSM row 107
$ws_BSM.Range("H107").Value = 2329.9707
$ws_BSM.Range("I107").Value = 1279.9584
$ws_BSM.Range("K107").Value = 1279.9584
$ws_BSM.Range("M107").Value = 640.0416

# BSM row 134
$ws_BSM.Range("H134").Value = 2213.1614
$ws_BSM.Range("I134").Value = 2213.1614
$ws_BSM.Range("K134").Value = 6639.4842
$ws_BSM.Range("M134").Value = -4104.4842

# CRP row 31
$ws_CRP.Range("H31").Value = 7333.6
$ws_CRP.Range("I31").Value = 8952.75
$ws_CRP.Range("K31").Value = 8952.75
$ws_CRP.Range("M31").Value = -8657.75

# CRP row 34
$ws_CRP.Range("H34").Value = 7333.6
$ws_CRP.Range("I34").Value = 8952.75
$ws_CRP.Range("K34").Value = 8952.75
$ws_CRP.Range("M34").Value = -8750.75

# CRP row 58
$ws_CRP.Range("H58").Value = 45259.13
$ws_CRP.Range("I58").Value = 45259.13
$ws_CRP.Range("K58").Value = 45259.13
$ws_CRP.Range("M58").Value = -45056.13

# CRP row 68
$ws_CRP.Range("H68").Value = 54250
$ws_CRP.Range("J68").Value = 54250
$ws_CRP.Range("L68").Value = 54250
$ws_CRP.Range("N68").Value = -55748

# CRP row 71
$ws_CRP.Range("H71").Value = 54250
$ws_CRP.Range("J71").Value = 54250
$ws_CRP.Range("L71").Value = 162750
$ws_CRP.Range("N71").Value = -170238

# CRP row 132
$ws_CRP.Range("H132").Value = 2341.9167
$ws_CRP.Range("I132").Value = 2191.182
$ws_CRP.Range("K132").Value = 6573.545999999999
$ws_CRP.Range("M132").Value = -4043.545999999999

# CRP row 134
$ws_CRP.Range("H134").Value = 56451.2
$ws_CRP.Range("I134").Value = 78000.64
$ws_CRP.Range("J134").Value = 6169.1665
$ws_CRP.Range("K134").Value = 234001.92
$ws_CRP.Range("L134").Value = 18507.4995
$ws_CRP.Range("M134").Value = -231466.92
$ws_CRP.Range("N134").Value = -23577.4995

# CRP row 136
$ws_CRP.Range("H136").Value = 45259.13
$ws_CRP.Range("I136").Value = 45259.13
$ws_CRP.Range("K136").Value = 135777.39
$ws_CRP.Range("M136").Value = -133227.39

# CUL row 2
$ws_CUL.Range("H2").Value = 46.785713
$ws_CUL.Range("I2").Value = 45
$ws_CUL.Range("K2").Value = 270
$ws_CUL.Range("M2").Value = -157

# CUL row 5
$ws_CUL.Range("H5").Value = 402.8
$ws_CUL.Range("I5").Value = 402.8
$ws_CUL.Range("J5").Value = 0
$ws_CUL.Range("K5").Value = 1208.4
$ws_CUL.Range("L5").Value = 0
$ws_CUL.Range("M5").Value = -1096.4
$ws_CUL.Range("N5").Value = ""

# CUL row 38
$ws_CUL.Range("H38").Value = 228.6875
$ws_CUL.Range("I38").Value = 353.75
$ws_CUL.Range("J38").Value = 187
$ws_CUL.Range("K38").Value = 1061.25
$ws_CUL.Range("L38").Value = 561
$ws_CUL.Range("M38").Value = -714.25
$ws_CUL.Range("N38").Value = -1255

# CUL row 122
$ws_CUL.Range("H122").Value = 395.75
$ws_CUL.Range("I122").Value = 395.75
$ws_CUL.Range("J122").Value = 0
$ws_CUL.Range("K122").Value = 3561.75
$ws_CUL.Range("L122").Value = 0
$ws_CUL.Range("M122").Value = -1111.75
$ws_CUL.Range("N122").Value = ""

# CUL row 123
$ws_CUL.Range("H123").Value = 23803.8
$ws_CUL.Range("J123").Value = 3000
$ws_CUL.Range("L123").Value = 9000
$ws_CUL.Range("N123").Value = -13900

# CUL row 132
$ws_CUL.Range("H132").Value = 2208.7778
$ws_CUL.Range("I132").Value = 1298.5
$ws_CUL.Range("J132").Value = 2937
$ws_CUL.Range("K132").Value = 11686.5
$ws_CUL.Range("L132").Value = 26433
$ws_CUL.Range("M132").Value = -9156.5
$ws_CUL.Range("N132").Value = -31493

# CUL row 135
$ws_CUL.Range("H135").Value = 402.8
$ws_CUL.Range("I135").Value = 402.8
$ws_CUL.Range("J135").Value = 0
$ws_CUL.Range("K135").Value = 3625.2
$ws_CUL.Range("L135").Value = 0
$ws_CUL.Range("M135").Value = -1090.2
$ws_CUL.Range("N135").Value = ""

# GSM row 132
$ws_GSM.Range("H132").Value = 54830.105
$ws_GSM.Range("I132").Value = 54830.105
$ws_GSM.Range("J132").Value = 0
$ws_GSM.Range("K132").Value = 164490.315
$ws_GSM.Range("L132").Value = 0
$ws_GSM.Range("M132").Value = -161960.315
$ws_GSM.Range("N132").Value = ""

# LTW row 7
$ws_LTW.Range("H7").Value = 4866.8335
$ws_LTW.Range("I7").Value = 2300.75
$ws_LTW.Range("K7").Value = 2300.75
$ws_LTW.Range("M7").Value = -2188.75

# LTW row 22
$ws_LTW.Range("H22").Value = 45716.24
$ws_LTW.Range("I22").Value = 93226.414
$ws_LTW.Range("K22").Value = 93226.414
$ws_LTW.Range("M22").Value = -92931.414

# LTW row 27
$ws_LTW.Range("H27").Value = 45716.24
$ws_LTW.Range("I27").Value = 93226.414
$ws_LTW.Range("K27").Value = 93226.414
$ws_LTW.Range("M27").Value = -93119.414

# LTW row 55
$ws_LTW.Range("H55").Value = 114.75
$ws_LTW.Range("I55").Value = 120.333336
$ws_LTW.Range("K55").Value = 120.333336
$ws_LTW.Range("M55").Value = 52.666664

# LTW row 126
$ws_LTW.Range("H126").Value = 4866.8335
$ws_LTW.Range("I126").Value = 2300.75
$ws_LTW.Range("K126").Value = 6902.25
$ws_LTW.Range("M126").Value = -4432.25
